$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 300-374 per diff ---
$ws.Range("D300").Value = 44511
$ws.Range("M300").Value = 16
$ws.Range("N300").Value = 285000
$ws.Range("O300").Value = 290000
$ws.Range("P300").Value = 287500
$ws.Range("S300").Value = 639

$ws.Range("D301").Value = 44511
$ws.Range("M301").Value = 16
$ws.Range("N301").Value = 255000
$ws.Range("O301").Value = 260000
$ws.Range("P301").Value = 257500
$ws.Range("S301").Value = 572

$ws.Range("D302").Value = 44511
$ws.Range("L302").Value = "Segunda"
$ws.Range("M302").Value = 16
$ws.Range("N302").Value = 235000
$ws.Range("O302").Value = 240000
$ws.Range("P302").Value = 237500
$ws.Range("S302").Value = 528

$ws.Range("D303").Value = 44306
$ws.Range("L303").Value = "Especial"
$ws.Range("M303").Value = 26

$ws.Range("D304").Value = 44306
$ws.Range("L304").Value = "Primera"
$ws.Range("M304").Value = 20

$ws.Range("D305").Value = 44469
$ws.Range("K305").Value = "Packham's Triumph"
$ws.Range("L305").Value = "Especial"
$ws.Range("M305").Value = 10
$ws.Range("N305").Value = 285000
$ws.Range("O305").Value = 290000
$ws.Range("P305").Value = 287500
$ws.Range("S305").Value = 639

$ws.Range("D306").Value = 44469
$ws.Range("M306").Value = 20
$ws.Range("N306").Value = 255000
$ws.Range("O306").Value = 260000
$ws.Range("P306").Value = 257500
$ws.Range("R306").Value = "Región de O'Higgins"
$ws.Range("S306").Value = 572

$ws.Range("D307").Value = 44469
$ws.Range("M307").Value = 16
$ws.Range("N307").Value = 235000
$ws.Range("O307").Value = 240000
$ws.Range("P307").Value = 237500
$ws.Range("R307").Value = "Región de O'Higgins"
$ws.Range("S307").Value = 528

$ws.Range("D308").Value = 44215
$ws.Range("K308").Value = "Bartlett de verano"
$ws.Range("L308").Value = "Primera"
$ws.Range("M308").Value = 20
$ws.Range("N308").Value = 290000
$ws.Range("O308").Value = 300000
$ws.Range("P308").Value = 295000
$ws.Range("S308").Value = 656

$ws.Range("D309").Value = 44407
$ws.Range("M309").Value = 16
$ws.Range("N309").Value = 225000
$ws.Range("O309").Value = 230000
$ws.Range("P309").Value = 227500
$ws.Range("R309").Value = "Provincia de Curicó"
$ws.Range("S309").Value = 506

$ws.Range("D310").Value = 44407
$ws.Range("M310").Value = 20
$ws.Range("N310").Value = 195000
$ws.Range("O310").Value = 200000
$ws.Range("P310").Value = 197500
$ws.Range("R310").Value = "Provincia de Curicó"
$ws.Range("S310").Value = 439

$ws.Range("D311").Value = 44504
$ws.Range("K311").Value = "Packham's Triumph"
$ws.Range("M311").Value = 8
$ws.Range("N311").Value = 280000
$ws.Range("O311").Value = 285000
$ws.Range("P311").Value = 282500
$ws.Range("S311").Value = 628

$ws.Range("D312").Value = 44504
$ws.Range("K312").Value = "Packham's Triumph"
$ws.Range("M312").Value = 10
$ws.Range("N312").Value = 260000
$ws.Range("O312").Value = 265000
$ws.Range("P312").Value = 262500
$ws.Range("S312").Value = 583

$ws.Range("D313").Value = 44504
$ws.Range("L313").Value = "Segunda"
$ws.Range("N313").Value = 235000
$ws.Range("O313").Value = 240000
$ws.Range("P313").Value = 237500
$ws.Range("S313").Value = 528

$ws.Range("D314").Value = 44246
$ws.Range("K314").Value = "Bartlett de verano"
$ws.Range("L314").Value = "Especial"
$ws.Range("M314").Value = 20
$ws.Range("N314").Value = 235000
$ws.Range("O314").Value = 240000
$ws.Range("P314").Value = 237500
$ws.Range("S314").Value = 528

$ws.Range("D315").Value = 44246
$ws.Range("K315").Value = "Bartlett de verano"
$ws.Range("L315").Value = "Primera"
$ws.Range("M315").Value = 20
$ws.Range("N315").Value = 205000
$ws.Range("O315").Value = 210000
$ws.Range("P315").Value = 207500
$ws.Range("S315").Value = 461

$ws.Range("D316").Value = 44505
$ws.Range("L316").Value = "Especial"
$ws.Range("M316").Value = 10
$ws.Range("N316").Value = 280000
$ws.Range("O316").Value = 285000
$ws.Range("P316").Value = 282500
$ws.Range("S316").Value = 628

$ws.Range("D317").Value = 44505
$ws.Range("L317").Value = "Primera"
$ws.Range("M317").Value = 16
$ws.Range("N317").Value = 260000
$ws.Range("O317").Value = 265000
$ws.Range("P317").Value = 262500
$ws.Range("S317").Value = 583

$ws.Range("D318").Value = 44505
$ws.Range("K318").Value = "Packham's Triumph"
$ws.Range("L318").Value = "Segunda"
$ws.Range("N318").Value = 230000
$ws.Range("O318").Value = 235000
$ws.Range("P318").Value = 232500
$ws.Range("S318").Value = 517

$ws.Range("K319").Value = "Packham's Triumph"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 20
$ws.Range("N319").Value = 255000
$ws.Range("O319").Value = 260000
$ws.Range("P319").Value = 257500
$ws.Range("S319").Value = 572

$ws.Range("D320").Value = 44487
$ws.Range("L320").Value = "Segunda"
$ws.Range("M320").Value = 20
$ws.Range("N320").Value = 225000
$ws.Range("O320").Value = 230000
$ws.Range("P320").Value = 227500
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 506

$ws.Range("D321").Value = 44487
$ws.Range("K321").Value = "Winter Nelis"
$ws.Range("M321").Value = 10
$ws.Range("N321").Value = 255000
$ws.Range("O321").Value = 260000
$ws.Range("P321").Value = 257500
$ws.Range("R321").Value = "Región de O'Higgins"
$ws.Range("S321").Value = 572

$ws.Range("D322").Value = 44487
$ws.Range("K322").Value = "Winter Nelis"
$ws.Range("M322").Value = 10
$ws.Range("N322").Value = 235000
$ws.Range("O322").Value = 240000
$ws.Range("P322").Value = 237500
$ws.Range("R322").Value = "Región de O'Higgins"
$ws.Range("S322").Value = 528

$ws.Range("K323").Value = "Packham's Triumph"
$ws.Range("L323").Value = "Especial"
$ws.Range("M323").Value = 26
$ws.Range("N323").Value = 250000
$ws.Range("O323").Value = 260000
$ws.Range("P323").Value = 255000
$ws.Range("S323").Value = 567

$ws.Range("K324").Value = "Packham's Triumph"
$ws.Range("L324").Value = "Primera"
$ws.Range("M324").Value = 20
$ws.Range("N324").Value = 230000
$ws.Range("O324").Value = 240000
$ws.Range("P324").Value = 235000
$ws.Range("S324").Value = 522

$ws.Range("D325").Value = 44425
$ws.Range("L325").Value = "Segunda"
$ws.Range("M325").Value = 18
$ws.Range("N325").Value = 200000
$ws.Range("O325").Value = 210000
$ws.Range("P325").Value = 205000
$ws.Range("R325").Value = "Provincia de Curicó"
$ws.Range("S325").Value = 456

$ws.Range("D326").Value = 44425
$ws.Range("K326").Value = "Winter Nelis"
$ws.Range("M326").Value = 24
$ws.Range("N326").Value = 230000
$ws.Range("O326").Value = 240000
$ws.Range("P326").Value = 235000
$ws.Range("R326").Value = "Provincia de Curicó"
$ws.Range("S326").Value = 522

$ws.Range("D327").Value = 44425
$ws.Range("K327").Value = "Winter Nelis"
$ws.Range("M327").Value = 18
$ws.Range("N327").Value = 210000
$ws.Range("O327").Value = 220000
$ws.Range("P327").Value = 215000
$ws.Range("R327").Value = "Provincia de Curicó"
$ws.Range("S327").Value = 478

$ws.Range("D328").Value = 44343
$ws.Range("N328").Value = 225000
$ws.Range("O328").Value = 230000
$ws.Range("P328").Value = 227500
$ws.Range("S328").Value = 506

$ws.Range("D329").Value = 44343
$ws.Range("N329").Value = 205000
$ws.Range("O329").Value = 210000
$ws.Range("P329").Value = 207500
$ws.Range("S329").Value = 461

$ws.Range("D330").Value = 44343
$ws.Range("N330").Value = 185000
$ws.Range("O330").Value = 190000
$ws.Range("P330").Value = 187500
$ws.Range("S330").Value = 417

$ws.Range("D331").Value = 44449
$ws.Range("N331").Value = 255000
$ws.Range("O331").Value = 260000
$ws.Range("P331").Value = 257500
$ws.Range("S331").Value = 572

$ws.Range("D332").Value = 44449
$ws.Range("N332").Value = 235000
$ws.Range("O332").Value = 240000
$ws.Range("P332").Value = 237500
$ws.Range("S332").Value = 528

$ws.Range("D333").Value = 44449
$ws.Range("N333").Value = 205000
$ws.Range("O333").Value = 210000
$ws.Range("P333").Value = 207500
$ws.Range("S333").Value = 461

$ws.Range("K334").Value = "Packham's Triumph"
$ws.Range("M334").Value = 20
$ws.Range("N334").Value = 235000
$ws.Range("O334").Value = 240000
$ws.Range("P334").Value = 237500
$ws.Range("S334").Value = 528

$ws.Range("K335").Value = "Packham's Triumph"

$ws.Range("K336").Value = "Packham's Triumph"
$ws.Range("M336").Value = 20

$ws.Range("D337").Value = 44392
$ws.Range("L337").Value = "Especial"
$ws.Range("M337").Value = 16
$ws.Range("N337").Value = 225000
$ws.Range("O337").Value = 230000
$ws.Range("P337").Value = 227500
$ws.Range("S337").Value = 506

$ws.Range("D338").Value = 44392
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 20
$ws.Range("N338").Value = 205000
$ws.Range("O338").Value = 210000
$ws.Range("P338").Value = 207500
$ws.Range("S338").Value = 461

$ws.Range("D339").Value = 44392
$ws.Range("K339").Value = "Winter Nelis"
$ws.Range("L339").Value = "Segunda"
$ws.Range("M339").Value = 16
$ws.Range("N339").Value = 185000
$ws.Range("O339").Value = 190000
$ws.Range("P339").Value = 187500
$ws.Range("S339").Value = 417

$ws.Range("D340").Value = 44286
$ws.Range("K340").Value = "Winter Nelis"
$ws.Range("L340").Value = "Primera"
$ws.Range("M340").Value = 22
$ws.Range("N340").Value = 195000
$ws.Range("O340").Value = 200000
$ws.Range("P340").Value = 197500
$ws.Range("S340").Value = 439

$ws.Range("D341").Value = 44286
$ws.Range("K341").Value = "Winter Nelis"
$ws.Range("L341").Value = "Segunda"
$ws.Range("M341").Value = 18
$ws.Range("N341").Value = 155000
$ws.Range("O341").Value = 160000
$ws.Range("P341").Value = 157500
$ws.Range("R341").Value = "Región de O'Higgins"
$ws.Range("S341").Value = 350

$ws.Range("D342").Value = 44473
$ws.Range("N342").Value = 245000
$ws.Range("O342").Value = 250000
$ws.Range("P342").Value = 247500
$ws.Range("R342").Value = "Región de O'Higgins"
$ws.Range("S342").Value = 550

$ws.Range("D343").Value = 44473
$ws.Range("M343").Value = 20
$ws.Range("N343").Value = 215000
$ws.Range("O343").Value = 220000
$ws.Range("P343").Value = 217500
$ws.Range("R343").Value = "Región de O'Higgins"
$ws.Range("S343").Value = 483

$ws.Range("K344").Value = "Packham's Triumph"
$ws.Range("R344").Value = "Provincia de Curicó"

$ws.Range("K345").Value = "Packham's Triumph"
$ws.Range("R345").Value = "Provincia de Curicó"

$ws.Range("K346").Value = "Packham's Triumph"
$ws.Range("M346").Value = 16
$ws.Range("N346").Value = 195000
$ws.Range("O346").Value = 200000
$ws.Range("P346").Value = 197500
$ws.Range("R346").Value = "Provincia de Curicó"
$ws.Range("S346").Value = 439

$ws.Range("D347").Value = 44400
$ws.Range("K347").Value = "Winter Nelis"
$ws.Range("L347").Value = "Especial"
$ws.Range("M347").Value = 16
$ws.Range("N347").Value = 235000
$ws.Range("O347").Value = 240000
$ws.Range("P347").Value = 237500
$ws.Range("S347").Value = 528

$ws.Range("D348").Value = 44400
$ws.Range("K348").Value = "Winter Nelis"
$ws.Range("L348").Value = "Primera"
$ws.Range("M348").Value = 20
$ws.Range("N348").Value = 215000
$ws.Range("O348").Value = 220000
$ws.Range("P348").Value = 217500
$ws.Range("S348").Value = 483

$ws.Range("D349").Value = 44400
$ws.Range("L349").Value = "Segunda"
$ws.Range("M349").Value = 20
$ws.Range("N349").Value = 175000
$ws.Range("O349").Value = 180000
$ws.Range("P349").Value = 177500
$ws.Range("S349").Value = 394

$ws.Range("K350").Value = "Packham's Triumph"
$ws.Range("L350").Value = "Primera"
$ws.Range("M350").Value = 20
$ws.Range("N350").Value = 250000
$ws.Range("O350").Value = 255000
$ws.Range("P350").Value = 252500
$ws.Range("S350").Value = 561

$ws.Range("D351").Value = 44484
$ws.Range("L351").Value = "Segunda"
$ws.Range("M351").Value = 16
$ws.Range("N351").Value = 220000
$ws.Range("O351").Value = 225000
$ws.Range("P351").Value = 222500
$ws.Range("S351").Value = 494

$ws.Range("D352").Value = 44484
$ws.Range("K352").Value = "Winter Nelis"
$ws.Range("N352").Value = 255000
$ws.Range("O352").Value = 260000
$ws.Range("P352").Value = 257500
$ws.Range("S352").Value = 572

$ws.Range("D353").Value = 44484
$ws.Range("K353").Value = "Winter Nelis"
$ws.Range("M353").Value = 10

$ws.Range("D354").Value = 44494
$ws.Range("L354").Value = "Especial"
$ws.Range("M354").Value = 10
$ws.Range("N354").Value = 295000
$ws.Range("O354").Value = 300000
$ws.Range("P354").Value = 297500
$ws.Range("S354").Value = 661

$ws.Range("D355").Value = 44494
$ws.Range("L355").Value = "Primera"
$ws.Range("M355").Value = 10
$ws.Range("N355").Value = 275000
$ws.Range("O355").Value = 280000
$ws.Range("P355").Value = 277500
$ws.Range("S355").Value = 617

$ws.Range("D356").Value = 44494
$ws.Range("L356").Value = "Segunda"
$ws.Range("M356").Value = 16
$ws.Range("N356").Value = 235000
$ws.Range("O356").Value = 240000
$ws.Range("P356").Value = 237500
$ws.Range("S356").Value = 528

$ws.Range("D357").Value = 44445
$ws.Range("L357").Value = "Primera"
$ws.Range("N357").Value = 235000
$ws.Range("O357").Value = 240000
$ws.Range("P357").Value = 237500
$ws.Range("S357").Value = 528

$ws.Range("D358").Value = 44445
$ws.Range("K358").Value = "Packham's Triumph"
$ws.Range("L358").Value = "Segunda"
$ws.Range("N358").Value = 205000
$ws.Range("O358").Value = 210000
$ws.Range("P358").Value = 207500
$ws.Range("S358").Value = 461

$ws.Range("K359").Value = "Packham's Triumph"
$ws.Range("L359").Value = "Primera"
$ws.Range("N359").Value = 225000
$ws.Range("O359").Value = 230000
$ws.Range("P359").Value = 227500
$ws.Range("S359").Value = 506

$ws.Range("D360").Value = 44301
$ws.Range("K360").Value = "Packham's Triumph"
$ws.Range("L360").Value = "Segunda"
$ws.Range("N360").Value = 205000
$ws.Range("O360").Value = 210000
$ws.Range("P360").Value = 207500
$ws.Range("S360").Value = 461

$ws.Range("D361").Value = 44301
$ws.Range("K361").Value = "Winter Nelis"
$ws.Range("L361").Value = "Primera"
$ws.Range("M361").Value = 20
$ws.Range("R361").Value = "Región de O'Higgins"

$ws.Range("D362").Value = 44301
$ws.Range("K362").Value = "Winter Nelis"
$ws.Range("L362").Value = "Segunda"
$ws.Range("R362").Value = "Región de O'Higgins"

$ws.Range("D363").Value = 44330
$ws.Range("K363").Value = "Winter Nelis"
$ws.Range("L363").Value = "Primera"
$ws.Range("N363").Value = 195000
$ws.Range("O363").Value = 200000
$ws.Range("P363").Value = 197500
$ws.Range("R363").Value = "Región de O'Higgins"
$ws.Range("S363").Value = 439

$ws.Range("D364").Value = 44295
$ws.Range("K364").Value = "Packham's Triumph"
$ws.Range("L364").Value = "Especial"
$ws.Range("M364").Value = 16
$ws.Range("N364").Value = 235000
$ws.Range("O364").Value = 240000
$ws.Range("P364").Value = 237500
$ws.Range("R364").Value = "Provincia de Curicó"
$ws.Range("S364").Value = 528

$ws.Range("D365").Value = 44295
$ws.Range("K365").Value = "Packham's Triumph"
$ws.Range("L365").Value = "Primera"
$ws.Range("N365").Value = 215000
$ws.Range("O365").Value = 220000
$ws.Range("P365").Value = 217500
$ws.Range("R365").Value = "Provincia de Curicó"
$ws.Range("S365").Value = 483

$ws.Range("D366").Value = 44295
$ws.Range("L366").Value = "Segunda"
$ws.Range("M366").Value = 20
$ws.Range("N366").Value = 165000
$ws.Range("O366").Value = 170000
$ws.Range("P366").Value = 167500
$ws.Range("R366").Value = "Provincia de Curicó"
$ws.Range("S366").Value = 372

$ws.Range("D367").Value = 44217
$ws.Range("K367").Value = "Bartlett de verano"
$ws.Range("N367").Value = 295000
$ws.Range("O367").Value = 300000
$ws.Range("P367").Value = 297500
$ws.Range("S367").Value = 661

$ws.Range("D368").Value = 44217
$ws.Range("K368").Value = "Bartlett de verano"
$ws.Range("M368").Value = 20

$ws.Range("D369").Value = 44509
$ws.Range("M369").Value = 24
$ws.Range("N369").Value = 285000
$ws.Range("O369").Value = 290000
$ws.Range("P369").Value = 287500
$ws.Range("R369").Value = "Región de O'Higgins"
$ws.Range("S369").Value = 639

$ws.Range("D370").Value = 44509
$ws.Range("N370").Value = 255000
$ws.Range("O370").Value = 260000
$ws.Range("P370").Value = 257500
$ws.Range("R370").Value = "Región de O'Higgins"
$ws.Range("S370").Value = 572

$ws.Range("D371").Value = 44509
$ws.Range("M371").Value = 18
$ws.Range("N371").Value = 235000
$ws.Range("O371").Value = 240000
$ws.Range("P371").Value = 237500
$ws.Range("R371").Value = "Región de O'Higgins"
$ws.Range("S371").Value = 528

$ws.Range("K372").Value = "Packham's Triumph"
$ws.Range("L372").Value = "Especial"
$ws.Range("M372").Value = 16
$ws.Range("N372").Value = 255000
$ws.Range("O372").Value = 260000
$ws.Range("P372").Value = 257500
$ws.Range("S372").Value = 572

$ws.Range("D373").Value = 44421
$ws.Range("N373").Value = 235000
$ws.Range("O373").Value = 240000
$ws.Range("P373").Value = 237500
$ws.Range("R373").Value = "Provincia de Curicó"
$ws.Range("S373").Value = 528

$ws.Range("D374").Value = 44421
$ws.Range("R374").Value = "Provincia de Curicó"

# --- Append new rows 375-377 ---
$ws.Range("A375").Value = 8
$ws.Range("B375").Value = "Terminal La Palmera de La Serena"
$ws.Range("C375").Value = "Coquimbo"
$ws.Range("D375").Value = 44421
$ws.Range("D375").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E375").Value = 4
$ws.Range("F375").Value = "Fruta"
$ws.Range("G375").Value = 100104
$ws.Range("H375").Value = "Frutos de pepita"
$ws.Range("I375").Value = 100104005
$ws.Range("J375").Value = "Pera"
$ws.Range("K375").Value = "Winter Nelis"
$ws.Range("L375").Value = "Primera"
$ws.Range("M375").Value = 10
$ws.Range("N375").Value = 235000
$ws.Range("O375").Value = 240000
$ws.Range("P375").Value = 237500
$ws.Range("Q375").Value = "`$/bins (450 kilos)"
$ws.Range("R375").Value = "Provincia de Curicó"
$ws.Range("S375").Value = 528
$ws.Range("T375").Value = 450

$ws.Range("A376").Value = 8
$ws.Range("B376").Value = "Terminal La Palmera de La Serena"
$ws.Range("C376").Value = "Coquimbo"
$ws.Range("D376").Value = 44302
$ws.Range("D376").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E376").Value = 4
$ws.Range("F376").Value = "Fruta"
$ws.Range("G376").Value = 100104
$ws.Range("H376").Value = "Frutos de pepita"
$ws.Range("I376").Value = 100104005
$ws.Range("J376").Value = "Pera"
$ws.Range("K376").Value = "Packham's Triumph"
$ws.Range("L376").Value = "Primera"
$ws.Range("M376").Value = 20
$ws.Range("N376").Value = 225000
$ws.Range("O376").Value = 230000
$ws.Range("P376").Value = 227500
$ws.Range("Q376").Value = "`$/bins (450 kilos)"
$ws.Range("R376").Value = "Región de O'Higgins"
$ws.Range("S376").Value = 506
$ws.Range("T376").Value = 450

$ws.Range("A377").Value = 8
$ws.Range("B377").Value = "Terminal La Palmera de La Serena"
$ws.Range("C377").Value = "Coquimbo"
$ws.Range("D377").Value = 44302
$ws.Range("D377").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E377").Value = 4
$ws.Range("F377").Value = "Fruta"
$ws.Range("G377").Value = 100104
$ws.Range("H377").Value = "Frutos de pepita"
$ws.Range("I377").Value = 100104005
$ws.Range("J377").Value = "Pera"
$ws.Range("K377").Value = "Packham's Triumph"
$ws.Range("L377").Value = "Segunda"
$ws.Range("M377").Value = 20
$ws.Range("N377").Value = 205000
$ws.Range("O377").Value = 210000
$ws.Range("P377").Value = 207500
$ws.Range("Q377").Value = "`$/bins (450 kilos)"
$ws.Range("R377").Value = "Región de O'Higgins"
$ws.Range("S377").Value = 461
$ws.Range("T377").Value = 450
